$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.981.78'
$ws.Range("E2").Value = '  +3.97%  '

$ws.Range("D3").Value = '3.241.86'
$ws.Range("E3").Value = '  +2.48%  '

$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").Value = '542.89'
$ws.Range("E5").Value = '  +2.62%  '

$ws.Range("D6").Value = '147.60'
$ws.Range("E6").Value = '  +5.46%  '

$ws.Range("E7").Value = '  -0.15%  '

$ws.Range("D8").Value = '0.528'
$ws.Range("E8").Value = '  -1.46%  '

$ws.Range("E9").Value = '  +0.99%  '

$ws.Range("E10").Value = '  +2.83%  '

$ws.Range("D11").Value = '0.437'
$ws.Range("E11").Value = '  -0.35%  '

$ws.Range("D12").Value = '3.792.02'
$ws.Range("E12").Value = '  +2.21%  '

$ws.Range("E13").Value = '  -1.89%  '

$ws.Range("D14").Value = '26.23'
$ws.Range("E14").Value = '  +2.02%  '

$ws.Range("E15").Value = '  +2.72%  '

$ws.Range("D16").Value = '60.896.95'
$ws.Range("E16").Value = '  +3.74%  '

$ws.Range("D17").Value = '3.244.95'
$ws.Range("E17").Value = '  +2.37%  '

$ws.Range("D18").Value = '6.32'
$ws.Range("E18").Value = '  +1.53%  '

$ws.Range("D19").Value = '13.36'
$ws.Range("E19").Value = '  +3.17%  '

$ws.Range("D20").Value = '8.37'
$ws.Range("E20").Value = '  +3.24%  '

$ws.Range("D21").Value = '378.02'
$ws.Range("E21").Value = '  +0.58%  '

$ws.Range("E23").Value = '  -0.08%  '

$ws.Range("D24").Value = '70.03'
$ws.Range("E24").Value = '  +0.49%  '

$ws.Range("E25").Value = '  +2.42%  '

$ws.Range("D26").Value = '8.67'
$ws.Range("E26").Value = '  +4.28%  '

$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.03%  '

$ws.Range("D28").Value = '0.0₃0914'
$ws.Range("E28").Value = '  +6.44%  '

$ws.Range("D29").Value = '22.60'
$ws.Range("E29").Value = '  +0.53%  '

$ws.Range("E30").Value = '  +2.06%  '

$ws.Range("E31").Value = '  +3.56%  '

$ws.Range("D32").Value = '5.38'
$ws.Range("E32").Value = '  +4.77%  '

$ws.Range("E33").Value = '  +6.98%  '

$ws.Range("D34").Value = '6.62'
$ws.Range("E34").Value = '  +5.20%  '

$ws.Range("D35").Value = '158.28'
$ws.Range("E35").Value = '  +0.99%  '

$ws.Range("E36").Value = '  +7.06%  '

$ws.Range("D37").Value = '26.48'
$ws.Range("E37").Value = '  +5.90%  '

$ws.Range("D38").Value = '2.819.68'
$ws.Range("E38").Value = '  +4.64%  '

$ws.Range("D39").Value = '0.0715'
$ws.Range("E39").Value = '  +3.36%  '

$ws.Range("D40").Value = '0.0315'
$ws.Range("E40").Value = '  +8.54%  '

$ws.Range("E41").Value = '  +2.20%  '

$ws.Range("E42").Value = '  -0.19%  '

$ws.Range("E43").Value = '  +2.06%  '

$ws.Range("D44").Value = '0.726'
$ws.Range("E44").Value = '  +0.69%  '

$ws.Range("B45").Value = 'RenzoRestakedETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D45").Value = '3.277.38'
$ws.Range("E45").Value = '  +2.21%  '

$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").Value = '0.105'
$ws.Range("E46").Value = '  +1.94%  '

$ws.Range("E47").Value = '  +2.69%  '

$ws.Range("D48").Value = '21.10'
$ws.Range("E48").Value = '  +5.40%  '

$ws.Range("D49").Value = '6.22'
$ws.Range("E49").Value = '  +0.22%  '

$ws.Range("D50").Value = '0.808'
$ws.Range("E50").Value = '  +8.24%  '

$ws.Range("E51").Value = '  -0.07%  '

